# New weekly price-report row for "Haba" (Vega Modelo de Temuco).
# Insert a fresh row at row 16 -- this pushes the existing data rows
# (old 16..72) down by one (new 17..73), exactly like pasting a new
# week's record at the top of the dated list while keeping history below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16 (shifts 16..72 -> 17..73).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the latest market record.
$ws.Cells.Item(16, 1).Value  = 10                        # Mercado ID
$ws.Cells.Item(16, 2).Value  = "Vega Modelo de Temuco"    # Mercado
$ws.Cells.Item(16, 3).Value  = "La Araucanía"             # Región
$ws.Cells.Item(16, 4).Value  = 44802                      # Fecha
$ws.Cells.Item(16, 5).Value  = 9                          # Codreg
$ws.Cells.Item(16, 6).Value  = 100112026                  # Categoría ID
$ws.Cells.Item(16, 7).Value  = "Haba"                      # Categoría
$ws.Cells.Item(16, 8).Value  = "Sin especificar"          # Variedad
$ws.Cells.Item(16, 9).Value  = "Primera"                  # Calidad
$ws.Cells.Item(16, 10).Value = 55                         # Volumen
$ws.Cells.Item(16, 11).Value = 14000                      # Precio mínimo
$ws.Cells.Item(16, 12).Value = 14000                      # Precio máximo
$ws.Cells.Item(16, 13).Value = 14000                      # Precio promedio ponderado
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"          # Unidad de comercialización
$ws.Cells.Item(16, 15).Value = "Provincia de Limarí"      # Origen
$ws.Cells.Item(16, 16).Value = 560                        # Precio $/Kg
$ws.Cells.Item(16, 17).Value = 25                         # Kg o Unidades
$ws.Cells.Item(16, 18).Value = "Hortaliza"                 # Clasificación
